$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previous table contents (keeps the existing header style on A1/B1).
$ws.Range("A1:B12").ClearContents()

# The new first column (user_id) reuses the same header style as the
# existing header cells; give the new header cell (C1) that same style by
# copying the format from the adjacent "reason" header cell.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row: user_id | name | reason
$ws.Cells.Item(1, 1).Value = "user_id"
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "reason"

# Data rows (user_id, name, reason)
$data = @(
    @("U2020-0058", "MUHAMMAD ZIDNI", "There are 25 data"),
    @("U2021-0003", "ADAM HAIKAL", "There are 6 data"),
    @("U2021-0004", "ADIB M P", "Data doesnt exist"),
    @("U2021-0007", "AGUS NURYADI", "There are 25 data"),
    @("U2021-0175", "NURUL FITRI", "There are 41 data"),
    @("U2021-0188", "RAIHAN ABDAN", "There are 5 data"),
    @("U2022-0021", "ELSYE", "There are 62 data")
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $row++
}
